# Applies the "ניהול בקשות" edit: two new data rows on the first sheet
# (גיליון1), a widened column G, and an updated active-cell selection /
# scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("גיליון1")

# --- New row 6 -------------------------------------------------------
$ws.Range("A6").Value = 57698
$ws.Range("B6").Value = "intel"
$ws.Range("C6").Value = "נקלט"
$ws.Range("D6").Value = "מוטי"
$ws.Range("F6").Value = 55000
$ws.Range("G6").Value = "רווחה"
$ws.Range("H6").Value = "טומי"

# --- New row 7 -------------------------------------------------------
$ws.Range("A7").Value = 48726
$ws.Range("B7").Value = "amd"
$ws.Range("C7").Value = "נקלט"
$ws.Range("D7").Value = "יחזקל"
$ws.Range("F7").Value = 9999
$ws.Range("G7").Value = "בריאות"
$ws.Range("H7").Value = "טומי"

# --- Column width change ---------------------------------------------
# (25 "characters" renders as ColumnWidth ~24.17 in this font metric)
$ws.Columns.Item(7).ColumnWidth = 24.17

# --- View state: scroll position + active selection -------------------
$ws.Range("C8").Select()
$excel.ActiveWindow.ScrollColumn = 3
